$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enterprises density (per 1000 people)" row (row 11): Micro / SMEs values
# updated to more precise figures. Values are stored as text in the
# workbook, so force a text number format before assigning, then restore
# the cell's style so no new style is introduced.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2.73"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "0.97"
$ws.Range("C11").Style = "Normal"

# "Enterprises (% of total)" row (row 12): Micro / SMEs / MSMEs values
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "71.54"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "25.39"
$ws.Range("C12").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.93"
$ws.Range("D12").Style = "Normal"

# "Value added to the economy (% of total)" row (row 16): Micro / MSMEs values
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "19.93"
$ws.Range("B16").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "44.33"
$ws.Range("D16").Style = "Normal"
